$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ------------------------------------------------------------------
# 1. Prepare row 15 (currently a completely empty row) by copying the
#    cell formatting from row 14, so the new data we place there picks
#    up the same styles (date / time / text / databar formats) used by
#    every other log row.
# ------------------------------------------------------------------
$ws.Range("B14:G14").Copy() | Out-Null
$ws.Range("B15:G15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Shift the existing log entries (rows 4-14) down by one row so a
#    new entry can be inserted at the top (row 4). Work from the
#    bottom up so data is not overwritten before it is copied.
# ------------------------------------------------------------------
for ($r = 14; $r -ge 4; $r--) {
    $target = $r + 1
    $ws.Range("B$target").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$target").Value2 = $ws.Range("C$r").Value2
    $ws.Range("D$target").Value2 = $ws.Range("D$r").Value2
    $ws.Range("E$target").Value2 = $ws.Range("E$r").Value2
    $ws.Range("F$target").Value2 = $ws.Range("F$r").Value2
    $ws.Range("G$target").Value2 = $ws.Range("G$r").Value2
}

# ------------------------------------------------------------------
# 3. Write the new development-log entry into row 4. The date (B4) and
#    fix/status (D4 - "enterWordAndTile()") carry over unchanged; the
#    time and the perception / reflections / progress fields describe
#    the new work on validating words against the wordlist.
# ------------------------------------------------------------------
$ws.Range("C4").Value2 = 0.35138888888888886
$ws.Range("E4").Value2 = "Testing validation against wordlist.txt"
$ws.Range("F4").Value2 = "enterwordAndile() method Now verifes if first part of user string (the word), is in the directory. Although a lot of unit testing is needed."
$ws.Range("G4").Value2 = 0.57999999999999996

# ------------------------------------------------------------------
# 4. Extend the conditional formatting (expression rule + data bar)
#    down to the newly used row 15.
# ------------------------------------------------------------------
$fcs1 = $ws.Range("B4:F14").FormatConditions
$fcs1.Item(1).ModifyAppliesToRange($ws.Range("B4:F15")) | Out-Null

$fcs2 = $ws.Range("G4:G14").FormatConditions
$fcs2.Item(1).ModifyAppliesToRange($ws.Range("G4:G15")) | Out-Null

# ------------------------------------------------------------------
# 5. Update the view so the new entry at the top of the log is visible
#    and selected, matching the editor's on-screen state after typing
#    the new entry.
# ------------------------------------------------------------------
$ws.Range("E4").Select() | Out-Null
